# cryptos.xlsx data refresh
# Updates Price (D) / Volume(1h) (E) for existing rows, and
# shifts Coin/Link/Price/Volume for rows 34-51 (one coin dropped
# out of the ranked list, NEARProtocol entered at the bottom).
#
# Every Price cell in this sheet is stored as text (thousand
# separators use '.', which Excel would otherwise misparse, and
# trailing zeros are significant, e.g. '1.000'). A leading
# apostrophe forces the written value to stay text -- exactly
# like typing it into Excel by hand -- instead of letting COM
# silently auto-convert it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '26.265.59'
    'E2' = '  +2.85%  '
    'D3' = '1.719.30'
    'E3' = '  +3.15%  '
    'D4' = '0.9994'
    'E4' = '  +0.03%  '
    'D5' = '239.84'
    'E5' = '  +0.67%  '
    'D6' = '0.9999'
    'E6' = '  -0.01%  '
    'D7' = '0.4722'
    'E7' = '  -1.51%  '
    'D8' = '0.2623'
    'E8' = '  -0.38%  '
    'D9' = '0.06198'
    'E9' = '  +0.47%  '
    'D10' = '1.715.10'
    'E10' = '  +2.90%  '
    'D11' = '0.07060'
    'E11' = '  -0.52%  '
    'D12' = '15.47'
    'E12' = '  +4.39%  '
    'D13' = '0.5985'
    'E13' = '  +1.42%  '
    'D14' = '4.431'
    'E14' = '  +1.37%  '
    'D15' = '76.15'
    'E15' = '  +1.39%  '
    'E16' = '  +0.05%  '
    'E17' = '  -0.01%  '
    'D18' = '26.275.06'
    'E18' = '  +2.86%  '
    'D19' = '0.000006810'
    'E19' = '  +0.88%  '
    'D20' = '11.53'
    'E20' = '  +0.48%  '
    'D21' = '1.936.01'
    'E21' = '  +3.23%  '
    'D22' = '4.539'
    'E22' = '  +2.69%  '
    'D23' = '8.708'
    'E23' = '  -0.25%  '
    'D24' = '5.245'
    'E24' = '  -0.69%  '
    'D25' = '135.07'
    'E25' = '  -0.58%  '
    'D26' = '15.20'
    'E26' = '  +1.06%  '
    'D27' = '1.769'
    'E27' = '  +3.18%  '
    'E28' = '  +0.65%  '
    'D29' = '106.82'
    'E29' = '  +1.61%  '
    'D30' = '3.941'
    'E30' = '  -0.84%  '
    'D31' = '3.688'
    'E31' = '  +1.25%  '
    'D32' = '0.07788'
    'E32' = '  +0.57%  '
    'D33' = '0.04503'
    'E33' = '  +6.33%  '
    'B34' = 'HuobiToken'
    'C34' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D34' = '2.615'
    'E34' = '  +0.63%  '
    'B35' = 'ARBITRUM'
    'C35' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D35' = '0.9765'
    'E35' = '  +2.76%  '
    'B36' = 'ImmutableX'
    'C36' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D36' = '0.6203'
    'E36' = '  +1.36%  '
    'B37' = 'TrustWalletToken'
    'C37' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D37' = '0.9309'
    'E37' = '  +7.69%  '
    'B38' = 'Quant'
    'C38' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D38' = '115.16'
    'E38' = '  +18.54%  '
    'B39' = 'MXToken'
    'C39' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D39' = '2.436'
    'E39' = '  -6.08%  '
    'B40' = 'RenderToken'
    'C40' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D40' = '1.926'
    'E40' = '  +4.05%  '
    'B41' = 'PaxDollar'
    'C41' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D41' = '1.000'
    'E41' = '  +0.08%  '
    'B42' = 'FraxShare'
    'C42' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D42' = '5.615'
    'E42' = '  +15.68%  '
    'B43' = 'VeChain'
    'C43' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D43' = '0.01480'
    'E43' = '  +0.79%  '
    'B44' = 'TheSandbox'
    'C44' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D44' = '0.3826'
    'E44' = '  +1.57%  '
    'B45' = 'Algorand'
    'C45' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D45' = '0.1178'
    'E45' = '  +5.10%  '
    'B46' = 'Aptos'
    'C46' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D46' = '6.313'
    'E46' = '  +1.68%  '
    'B47' = 'Cronos'
    'C47' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D47' = '0.05263'
    'E47' = '  -0.07%  '
    'B48' = 'EnergySwap'
    'C48' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D48' = '7.789'
    'E48' = '  +5.58%  '
    'B49' = 'Elrond'
    'C49' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    'D49' = '30.41'
    'E49' = '  +2.16%  '
    'B50' = 'Decentraland'
    'C50' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'D50' = '0.3382'
    'E50' = '  +1.54%  '
    'B51' = 'NEARProtocol'
    'C51' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D51' = '1.219'
    'E51' = '  +1.87%  '
}

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    if ($ref.StartsWith('D')) {
        # Price column: keep as text (see note above)
        $ws.Range($ref).Value = "'" + $value
    } else {
        $ws.Range($ref).Value = $value
    }
}

